# Additional test case added. Rework of import classes.
$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Simplify the (stray multi-area) selections on the existing sheets down
#    to plain single-cell selections - same active cell, no more "C2:C5 .."
#    leading range glued onto the sqref.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Sheet1").Range("D1").Select()
$wb.Worksheets.Item("Sheet2").Range("E6").Select()
$wb.Worksheets.Item("Sheet3(empty)").Range("A1").Select()
$wb.Worksheets.Item("Sheet4").Range("C2").Select()
$wb.Worksheets.Item("Sheet5").Range("I9").Select()
$wb.Worksheets.Item("Sheet6").Range("A27").Select()
$wb.Worksheets.Item("Sheet7").Range("A22").Select()

# ---------------------------------------------------------------------------
# 2) Fix up the "general" number format (id 169) used on Sheet4 A2:A5/C2:C5 -
#    it should render like the neighbouring custom format (id 168) instead
#    of plain General.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("A2:A5").NumberFormat = '[$-415]ge\ner"al"'
$ws4.Range("C2:C5").NumberFormat = '[$-415]ge\ner"al"'

# ---------------------------------------------------------------------------
# 3) Split the old "testAccounts" sheet: keep its original data under a new
#    name ("Sheet9") is wrong -- instead insert a brand new small sheet
#    ("Sheet9") right before "testAccounts" and keep "testAccounts" (and its
#    data) as the last sheet, just re-created so it keeps sqref/selection in
#    the simplified single-area form too.
# ---------------------------------------------------------------------------
$old = $wb.Worksheets.Item("testAccounts")
$old.Delete()

$afterSheet7 = $wb.Worksheets.Item(7)
$sheet9 = $wb.Worksheets.Add($null, $afterSheet7)
$sheet9.Name = "Sheet9"

$sheet9.Range("B1").Value = "second"
$sheet9.Range("C1").Value = "third"
$sheet9.Range("A2").Value = "a"
$sheet9.Range("B2").Value = 45
$sheet9.Range("A3").Value = "s"
$sheet9.Range("B3").Value = 5
$sheet9.Range("A4").Value = "d"
$sheet9.Range("B4").Value = 67

$afterSheet9 = $wb.Worksheets.Item(8)
$ta = $wb.Worksheets.Add($null, $afterSheet9)
$ta.Name = "testAccounts"

$ta.Range("A1").Value = "user"
$ta.Range("B1").Value = "pass"
$ta.Range("C1").Value = "lic_exp"
$ta.Range("D1").Value = "uwagi"

$ta.Range("A2").Value = "test1"
$ta.Range("B2").Value = "test123"
$ta.Range("C2").Value = "zawsze aktualna"

$ta.Range("A3").Value = "test2"
$ta.Range("B3").Value = "test123"
$ta.Range("D3").Value = "konto zablokowane"

$ta.Range("A4").Value = "test3"
$ta.Range("B4").Value = "test123"
$ta.Range("C4").Value = 40816
$ta.Range("C4").NumberFormat = "[$-809]dd/mm/yyyy"

$ta.Range("F17").Select()

# ---------------------------------------------------------------------------
# 4) "Sheet9" is the sheet left active/selected in the final workbook.
# ---------------------------------------------------------------------------
$sheet9.Range("A1").Select()
